$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cells that actually changed in rows 2-10 ---
$ws.Range("G2").Value = "MJ2"
$ws.Range("H2").Value = "Good"

$ws.Range("G3").Value = "MJ2"
$ws.Range("H3").Value = "Good"

$ws.Range("G4").Value = "MJ2"
$ws.Range("H4").Value = "Good"

$ws.Range("G5").Value = "MJ2"
$ws.Range("H5").Value = "Good"

$ws.Range("G6").Value = "EA2"
$ws.Range("H6").Value = "Good"

$ws.Range("G7").Value = "EA2"
$ws.Range("H7").Value = "Good"

$ws.Range("F8").Value = "SOI TRIPURA"
$ws.Range("G8").Value = "EA2"
$ws.Range("H8").Value = "Good"

$ws.Range("G9").Value = "EA2"
$ws.Range("H9").Value = "Good"

$ws.Range("F10").Value = "SOI ASSAM"
$ws.Range("G10").Value = "MO1"
$ws.Range("H10").Value = "Good"

# --- Row 11 previously held temp placeholder data ("as"/"asd"/"sd"/"ss"); replace with real row ---
$ws.Range("A11").Value = "Electronics"
$ws.Range("B11").Value = "Name10"
$ws.Range("C11").Value = "Make10"
$ws.Range("D11").Value = "Model10"
$ws.Range("E11").Value = "Serial12"
$ws.Range("F11").Value = "SOI ASSAM"
$ws.Range("G11").Value = "MO1"
$ws.Range("H11").Value = "Good"

# --- New rows 12-17 ---
$ws.Range("A12").Value = "Drone Equipment"
$ws.Range("B12").Value = "Name11"
$ws.Range("C12").Value = "Make11"
$ws.Range("D12").Value = "Model11"
$ws.Range("E12").Value = "Serial13"
$ws.Range("F12").Value = "SOI ASSAM"
$ws.Range("G12").Value = "MO1"
$ws.Range("H12").Value = "Good"

$ws.Range("A13").Value = "Electronics"
$ws.Range("B13").Value = "Name12"
$ws.Range("C13").Value = "Make12"
$ws.Range("D13").Value = "Model12"
$ws.Range("E13").Value = "Serial14"
$ws.Range("F13").Value = "SOI ASSAM"
$ws.Range("G13").Value = "MO1"
$ws.Range("H13").Value = "Good"

$ws.Range("A14").Value = "DGPS Equipment"
$ws.Range("B14").Value = "Name13"
$ws.Range("C14").Value = "Make13"
$ws.Range("D14").Value = "Model13"
$ws.Range("E14").Value = "Serial15"
$ws.Range("F14").Value = "SOI ASSAM"
$ws.Range("G14").Value = "ES1"
$ws.Range("H14").Value = "Good"

$ws.Range("A15").Value = "DGPS Equipment"
$ws.Range("B15").Value = "Name14"
$ws.Range("C15").Value = "Make14"
$ws.Range("D15").Value = "Model14"
$ws.Range("E15").Value = "Serial16"
$ws.Range("F15").Value = "SOI ASSAM"
$ws.Range("G15").Value = "ES1"
$ws.Range("H15").Value = "Good"

$ws.Range("A16").Value = "Drone Equipment"
$ws.Range("B16").Value = "Name15"
$ws.Range("C16").Value = "Make15"
$ws.Range("D16").Value = "Model15"
$ws.Range("E16").Value = "Serial17"
$ws.Range("F16").Value = "SOI ASSAM"
$ws.Range("G16").Value = "ES1"
$ws.Range("H16").Value = "Good"

$ws.Range("A17").Value = "Drone Equipment"
$ws.Range("B17").Value = "Name16"
$ws.Range("C17").Value = "Make16"
$ws.Range("D17").Value = "Model16"
$ws.Range("E17").Value = "Serial18"
$ws.Range("F17").Value = "SOI ASSAM"
$ws.Range("G17").Value = "ES1"
$ws.Range("H17").Value = "Good"

# --- Columns that Excel auto-sized ("best fit") after the new rows were typed in ---
$ws.Columns.Item(1).ColumnWidth = 15.833333333333334
$ws.Columns.Item(6).ColumnWidth = 11.166666666666666
$ws.Columns.Item(8).ColumnWidth = 8.833333333333334
$ws.Columns.Item(9).ColumnWidth = 12.333333333333334

# --- Selection left where the user last clicked ---
$ws.Range("P3").Select()

$wb.Save()
